# Updates the "cryptos" price/volume table (rows 2-51) on Sheet1 to match the
# latest scrape: refreshed Price (col D) / Volume(1h) (col E) figures, plus a
# few rows whose ranking swapped with a neighbor (Coin name + Link also change).
#
# Note: several Price values are plain decimals (e.g. "0.602", "41.73") which
# Excel's COM layer would otherwise auto-convert to a numeric cell. They are
# written with a leading apostrophe (the classic Excel "force text" prefix) to
# keep them as text like the rest of the column, and the cell Style is reset
# back to "Normal" afterwards so no stray text-number-format style is left on
# the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.306.19'
$ws.Range('E2').Value = '  -1.56%  '
$ws.Range('D3').Value = '2.146.96'
$ws.Range('E3').Value = '  -2.96%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'236.26"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.90%  '
$ws.Range('D6').Value = "'0.602"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.73%  '
$ws.Range('D7').Value = "'70.03"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.51%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = "'0.571"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.05%  '
$ws.Range('D10').Value = "'39.26"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.83%  '
$ws.Range('D11').Value = "'0.0896"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.65%  '
$ws.Range('D12').Value = "'53.62"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.78%  '
$ws.Range('D13').Value = "'0.0994"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.32%  '
$ws.Range('D14').Value = "'6.62"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.97%  '
$ws.Range('D15').Value = '2.462.94'
$ws.Range('E15').Value = '  -3.21%  '
$ws.Range('D16').Value = "'14.17"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').Value = '2.156.23'
$ws.Range('E17').Value = '  -2.28%  '
$ws.Range('D18').Value = "'0.779"
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = '41.068.01'
$ws.Range('E19').Value = '  -1.80%  '
$ws.Range('E20').Value = '  -4.73%  '
$ws.Range('D21').Value = "'68.92"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.90%  '
$ws.Range('E22').Value = '  -6.78%  '
$ws.Range('D23').Value = "'9.60"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -10.16%  '
$ws.Range('D24').Value = "'225.21"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.66%  '
$ws.Range('D25').Value = "'1.95"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.96%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = "'10.58"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.55%  '
$ws.Range('D28').Value = "'3.32"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -8.83%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E29').Value = '  -4.34%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = "'2.16"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.77%  '
$ws.Range('D31').Value = "'170.96"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.18%  '
$ws.Range('D32').Value = "'19.65"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.80%  '
$ws.Range('D33').Value = "'31.41"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.75%  '
$ws.Range('D34').Value = "'0.0755"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.76%  '
$ws.Range('E35').Value = '  -9.32%  '
$ws.Range('E36').Value = '  -3.62%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = "'4.23"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = "'0.102"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.54%  '
$ws.Range('D39').Value = "'0.0291"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.76%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = "'11.73"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -14.80%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').Value = "'2.04"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.82%  '
$ws.Range('E42').Value = '  -7.04%  '
$ws.Range('D43').Value = "'57.85"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -10.28%  '
$ws.Range('D44').Value = "'0.187"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.61%  '
$ws.Range('D45').Value = "'8.20"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.17%  '
$ws.Range('D46').Value = "'0.0958"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.53%  '
$ws.Range('D47').Value = "'97.30"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.29%  '
$ws.Range('E48').Value = '  -3.17%  '
$ws.Range('E49').Value = '  -5.18%  '
$ws.Range('E50').Value = '  -3.13%  '
$ws.Range('E51').Value = '  -8.35%  '
